$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.722831
$ws.Range("N2").Value = 11.168493
$ws.Range("O2").Value = 0.2042994277655142
$ws.Range("P2").Value = 0.2042994277655142
$ws.Range("Q2").Value = 151.7179886108647
$ws.Range("R2").Value = 1365.461897497782
$ws.Range("S2").Value = 0.004344502168738791
$ws.Range("T2").Value = 0.004344502168738791

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.746044
$ws.Range("N3").Value = 5.238131999999999
$ws.Range("O3").Value = 0.09581842153280916
$ws.Range("P3").Value = 0.09581842153280916
$ws.Range("Q3").Value = 71.15721441721865
$ws.Range("R3").Value = 640.4149297549679
$ws.Range("S3").Value = 0.002037613833320222
$ws.Range("T3").Value = 0.002037613833320222

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.75339133333333
$ws.Range("H4").Value = 122.260174
$ws.Range("I4").Value = 0.02126536631186857
$ws.Range("J4").Value = 0.02126536631186857
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 11.46242466666667
$ws.Range("N4").Value = 34.387274
$ws.Range("O4").Value = 0.6290285001401661
$ws.Range("P4").Value = 0.6290285001401661
$ws.Range("Q4").Value = 467.1326780695195
$ws.Range("R4").Value = 4204.194102625676
$ws.Range("S4").Value = 0.0133765214760859
$ws.Range("T4").Value = 0.0133765214760859

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.75339133333333
$ws.Range("H5").Value = 122.260174
$ws.Range("I5").Value = 0.02126536631186857
$ws.Range("J5").Value = 0.02126536631186857
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.291125333333333
$ws.Range("N5").Value = 3.873376
$ws.Range("O5").Value = 0.07085365056151052
$ws.Range("P5").Value = 0.07085365056151052
$ws.Range("Q5").Value = 52.61773596971377
$ws.Range("R5").Value = 473.559623727424
$ws.Range("S5").Value = 0.001506728833723653
$ws.Range("T5").Value = 0.001506728833723653

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1689.289306666667
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.8814813868902838
$ws.Range("J6").Value = 0.8814813868902838
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.722831
$ws.Range("N6").Value = 11.168493
$ws.Range("O6").Value = 0.2042994277655142
$ws.Range("P6").Value = 0.2042994277655142
$ws.Range("Q6").Value = 6288.938598827172
$ws.Range("R6").Value = 56600.44738944455
$ws.Range("S6").Value = 0.1800861429276368
$ws.Range("T6").Value = 0.1800861429276368

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1689.289306666667
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.8814813868902838
$ws.Range("J7").Value = 0.8814813868902838
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.746044
$ws.Range("N7").Value = 5.238131999999999
$ws.Range("O7").Value = 0.09581842153280916
$ws.Range("P7").Value = 0.09581842153280916
$ws.Range("Q7").Value = 2949.573458169493
$ws.Range("R7").Value = 26546.16112352543
$ws.Range("S7").Value = 0.08446215510237845
$ws.Range("T7").Value = 0.08446215510237845

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1689.289306666667
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.8814813868902838
$ws.Range("J8").Value = 0.8814813868902838
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 11.46242466666667
$ws.Range("N8").Value = 34.387274
$ws.Range("O8").Value = 0.6290285001401661
$ws.Range("P8").Value = 0.6290285001401661
$ws.Range("Q8").Value = 19363.35141787223
$ws.Range("R8").Value = 174270.1627608501
$ws.Range("S8").Value = 0.5544769146970686
$ws.Range("T8").Value = 0.5544769146970686

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1689.289306666667
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.8814813868902838
$ws.Range("J9").Value = 0.8814813868902838
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.291125333333333
$ws.Range("N9").Value = 3.873376
$ws.Range("O9").Value = 0.07085365056151052
$ws.Range("P9").Value = 0.07085365056151052
$ws.Range("Q9").Value = 2181.084219166435
$ws.Range("R9").Value = 19629.75797249792
$ws.Range("S9").Value = 0.06245617416319983
$ws.Range("T9").Value = 0.06245617416319983

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 100.9654023333333
$ws.Range("H10").Value = 302.896207
$ws.Range("I10").Value = 0.05268435816499466
$ws.Range("J10").Value = 0.05268435816499466
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.722831
$ws.Range("N10").Value = 11.168493
$ws.Range("O10").Value = 0.2042994277655142
$ws.Range("P10").Value = 0.2042994277655142
$ws.Range("Q10").Value = 375.8771297340056
$ws.Range("R10").Value = 3382.894167606051
$ws.Range("S10").Value = 0.0107633842253018
$ws.Range("T10").Value = 0.0107633842253018

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 100.9654023333333
$ws.Range("H11").Value = 302.896207
$ws.Range("I11").Value = 0.05268435816499466
$ws.Range("J11").Value = 0.05268435816499466
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.746044
$ws.Range("N11").Value = 5.238131999999999
$ws.Range("O11").Value = 0.09581842153280916
$ws.Range("P11").Value = 0.09581842153280916
$ws.Range("Q11").Value = 176.2900349517026
$ws.Range("R11").Value = 1586.610314565324
$ws.Range("S11").Value = 0.005048132038838954
$ws.Range("T11").Value = 0.005048132038838954

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 100.9654023333333
$ws.Range("H12").Value = 302.896207
$ws.Range("I12").Value = 0.05268435816499466
$ws.Range("J12").Value = 0.05268435816499466
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 11.46242466666667
$ws.Range("N12").Value = 34.387274
$ws.Range("O12").Value = 0.6290285001401661
$ws.Range("P12").Value = 0.6290285001401661
$ws.Range("Q12").Value = 1157.308318185524
$ws.Range("R12").Value = 10415.77486366972
$ws.Range("S12").Value = 0.0331399627973739
$ws.Range("T12").Value = 0.0331399627973739

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 100.9654023333333
$ws.Range("H13").Value = 302.896207
$ws.Range("I13").Value = 0.05268435816499466
$ws.Range("J13").Value = 0.05268435816499466
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.291125333333333
$ws.Range("N13").Value = 3.873376
$ws.Range("O13").Value = 0.07085365056151052
$ws.Range("P13").Value = 0.07085365056151052
$ws.Range("Q13").Value = 130.3589887427591
$ws.Range("R13").Value = 1173.230898684832
$ws.Range("S13").Value = 0.003732879103479995
$ws.Range("T13").Value = 0.003732879103479995

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 85.41274733333334
$ws.Range("H14").Value = 256.238242
$ws.Range("I14").Value = 0.04456888863285297
$ws.Range("J14").Value = 0.04456888863285297
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.722831
$ws.Range("N14").Value = 11.168493
$ws.Range("O14").Value = 0.2042994277655142
$ws.Range("P14").Value = 0.2042994277655142
$ws.Range("Q14").Value = 317.9772235677007
$ws.Range("R14").Value = 2861.795012109306
$ws.Range("S14").Value = 0.009105398443836793
$ws.Range("T14").Value = 0.009105398443836793

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 85.41274733333334
$ws.Range("H15").Value = 256.238242
$ws.Range("I15").Value = 0.04456888863285297
$ws.Range("J15").Value = 0.04456888863285297
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.746044
$ws.Range("N15").Value = 5.238131999999999
$ws.Range("O15").Value = 0.09581842153280916
$ws.Range("P15").Value = 0.09581842153280916
$ws.Range("Q15").Value = 149.1344150048827
$ws.Range("R15").Value = 1342.209735043944
$ws.Range("S15").Value = 0.004270520558271532
$ws.Range("T15").Value = 0.004270520558271532

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 85.41274733333334
$ws.Range("H16").Value = 256.238242
$ws.Range("I16").Value = 0.04456888863285297
$ws.Range("J16").Value = 0.04456888863285297
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 11.46242466666667
$ws.Range("N16").Value = 34.387274
$ws.Range("O16").Value = 0.6290285001401661
$ws.Range("P16").Value = 0.6290285001401661
$ws.Range("Q16").Value = 979.0371818813676
$ws.Range("R16").Value = 8811.334636932308
$ws.Range("S16").Value = 0.0280351011696376
$ws.Range("T16").Value = 0.0280351011696376

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 85.41274733333334
$ws.Range("H17").Value = 256.238242
$ws.Range("I17").Value = 0.04456888863285297
$ws.Range("J17").Value = 0.04456888863285297
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.291125333333333
$ws.Range("N17").Value = 3.873376
$ws.Range("O17").Value = 0.07085365056151052
$ws.Range("P17").Value = 0.07085365056151052
$ws.Range("Q17").Value = 110.2785618716658
$ws.Range("R17").Value = 992.507056844992
$ws.Range("S17").Value = 0.003157868461107043
$ws.Range("T17").Value = 0.003157868461107043
